$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, shifting existing rows 64-135 down to 65-136
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new price record
$ws.Range("A64").Value = 10
$ws.Range("B64").Value = "Vega Modelo de Temuco"
$ws.Range("C64").Value = "La Araucanía"
$ws.Range("D64").Value = 45271
$ws.Range("E64").Value = 9
$ws.Range("F64").Value = 100112022
$ws.Range("G64").Value = "Arveja Verde"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 35
$ws.Range("K64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("M64").Value = 25000
$ws.Range("N64").Value = "$/saco 25 kilos"
$ws.Range("O64").Value = "Región del Maule"
$ws.Range("P64").Value = 1000
$ws.Range("Q64").Value = 25
$ws.Range("R64").Value = "Hortaliza"
